$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dimension-relevant existing rows (176-229): shift price-history data down by
# one week (2 rows) and insert the newest week (date 44468) at rows 176-177.

# Row 176
$ws.Cells.Item(176, 4).Value = 44468
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 3120
$ws.Cells.Item(176, 11).Value = 450
$ws.Cells.Item(176, 12).Value = 500
$ws.Cells.Item(176, 13).Value = 475
$ws.Cells.Item(176, 16).Value = 238

# Row 177
$ws.Cells.Item(177, 4).Value = 44468
$ws.Cells.Item(177, 9).Value = "Segunda"
$ws.Cells.Item(177, 10).Value = 1360
$ws.Cells.Item(177, 11).Value = 350
$ws.Cells.Item(177, 12).Value = 400
$ws.Cells.Item(177, 13).Value = 375
$ws.Cells.Item(177, 16).Value = 188

# Row 178
$ws.Cells.Item(178, 4).Value = 44312
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 3000
$ws.Cells.Item(178, 11).Value = 550
$ws.Cells.Item(178, 12).Value = 600
$ws.Cells.Item(178, 13).Value = 575
$ws.Cells.Item(178, 16).Value = 288

# Row 179
$ws.Cells.Item(179, 4).Value = 44312
$ws.Cells.Item(179, 9).Value = "Segunda"
$ws.Cells.Item(179, 10).Value = 1500
$ws.Cells.Item(179, 11).Value = 450
$ws.Cells.Item(179, 12).Value = 500
$ws.Cells.Item(179, 13).Value = 475
$ws.Cells.Item(179, 16).Value = 238

# Row 180
$ws.Cells.Item(180, 4).Value = 44386
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 3460
$ws.Cells.Item(180, 11).Value = 450
$ws.Cells.Item(180, 12).Value = 500
$ws.Cells.Item(180, 13).Value = 475
$ws.Cells.Item(180, 16).Value = 238

# Row 181
$ws.Cells.Item(181, 4).Value = 44386
$ws.Cells.Item(181, 9).Value = "Segunda"
$ws.Cells.Item(181, 10).Value = 1600
$ws.Cells.Item(181, 11).Value = 350
$ws.Cells.Item(181, 12).Value = 400
$ws.Cells.Item(181, 13).Value = 375
$ws.Cells.Item(181, 16).Value = 188

# Row 182
$ws.Cells.Item(182, 4).Value = 44463
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 3400
$ws.Cells.Item(182, 11).Value = 450
$ws.Cells.Item(182, 12).Value = 500
$ws.Cells.Item(182, 13).Value = 475
$ws.Cells.Item(182, 16).Value = 238

# Row 183
$ws.Cells.Item(183, 4).Value = 44463
$ws.Cells.Item(183, 9).Value = "Segunda"
$ws.Cells.Item(183, 10).Value = 1500
$ws.Cells.Item(183, 11).Value = 350
$ws.Cells.Item(183, 12).Value = 400
$ws.Cells.Item(183, 13).Value = 375
$ws.Cells.Item(183, 16).Value = 188

# Row 184
$ws.Cells.Item(184, 4).Value = 44251
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 3600
$ws.Cells.Item(184, 11).Value = 450
$ws.Cells.Item(184, 12).Value = 500
$ws.Cells.Item(184, 13).Value = 475
$ws.Cells.Item(184, 16).Value = 238

# Row 185
$ws.Cells.Item(185, 4).Value = 44251
$ws.Cells.Item(185, 9).Value = "Segunda"
$ws.Cells.Item(185, 10).Value = 1760
$ws.Cells.Item(185, 11).Value = 350
$ws.Cells.Item(185, 12).Value = 400
$ws.Cells.Item(185, 13).Value = 375
$ws.Cells.Item(185, 16).Value = 188

# Row 186
$ws.Cells.Item(186, 4).Value = 44433
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 3500
$ws.Cells.Item(186, 11).Value = 450
$ws.Cells.Item(186, 12).Value = 500
$ws.Cells.Item(186, 13).Value = 475
$ws.Cells.Item(186, 16).Value = 238

# Row 187
$ws.Cells.Item(187, 4).Value = 44433
$ws.Cells.Item(187, 9).Value = "Segunda"
$ws.Cells.Item(187, 10).Value = 1600
$ws.Cells.Item(187, 11).Value = 350
$ws.Cells.Item(187, 12).Value = 400
$ws.Cells.Item(187, 13).Value = 375
$ws.Cells.Item(187, 16).Value = 188

# Row 188
$ws.Cells.Item(188, 4).Value = 44221
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 3000
$ws.Cells.Item(188, 11).Value = 450
$ws.Cells.Item(188, 12).Value = 500
$ws.Cells.Item(188, 13).Value = 475
$ws.Cells.Item(188, 16).Value = 238

# Row 189
$ws.Cells.Item(189, 4).Value = 44221
$ws.Cells.Item(189, 9).Value = "Segunda"
$ws.Cells.Item(189, 10).Value = 1600
$ws.Cells.Item(189, 11).Value = 350
$ws.Cells.Item(189, 12).Value = 400
$ws.Cells.Item(189, 13).Value = 375
$ws.Cells.Item(189, 16).Value = 188

# Row 190
$ws.Cells.Item(190, 4).Value = 44316
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 3400
$ws.Cells.Item(190, 11).Value = 500
$ws.Cells.Item(190, 12).Value = 600
$ws.Cells.Item(190, 13).Value = 550
$ws.Cells.Item(190, 16).Value = 275

# Row 191
$ws.Cells.Item(191, 4).Value = 44316
$ws.Cells.Item(191, 9).Value = "Segunda"
$ws.Cells.Item(191, 10).Value = 1600
$ws.Cells.Item(191, 11).Value = 400
$ws.Cells.Item(191, 12).Value = 450
$ws.Cells.Item(191, 13).Value = 425
$ws.Cells.Item(191, 16).Value = 212

# Row 192
$ws.Cells.Item(192, 4).Value = 44279
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 3400
$ws.Cells.Item(192, 11).Value = 450
$ws.Cells.Item(192, 12).Value = 500
$ws.Cells.Item(192, 13).Value = 475
$ws.Cells.Item(192, 16).Value = 238

# Row 193
$ws.Cells.Item(193, 4).Value = 44279
$ws.Cells.Item(193, 9).Value = "Segunda"
$ws.Cells.Item(193, 10).Value = 1600
$ws.Cells.Item(193, 11).Value = 350
$ws.Cells.Item(193, 12).Value = 400
$ws.Cells.Item(193, 13).Value = 375
$ws.Cells.Item(193, 16).Value = 188

# Row 194
$ws.Cells.Item(194, 4).Value = 44363
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 3440
$ws.Cells.Item(194, 11).Value = 500
$ws.Cells.Item(194, 12).Value = 600
$ws.Cells.Item(194, 13).Value = 550
$ws.Cells.Item(194, 16).Value = 275

# Row 195
$ws.Cells.Item(195, 4).Value = 44363
$ws.Cells.Item(195, 9).Value = "Segunda"
$ws.Cells.Item(195, 10).Value = 1600
$ws.Cells.Item(195, 11).Value = 400
$ws.Cells.Item(195, 12).Value = 450
$ws.Cells.Item(195, 13).Value = 425
$ws.Cells.Item(195, 16).Value = 212

# Row 196
$ws.Cells.Item(196, 4).Value = 44277
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 2800
$ws.Cells.Item(196, 11).Value = 450
$ws.Cells.Item(196, 12).Value = 500
$ws.Cells.Item(196, 13).Value = 475
$ws.Cells.Item(196, 16).Value = 238

# Row 197
$ws.Cells.Item(197, 4).Value = 44277
$ws.Cells.Item(197, 9).Value = "Segunda"
$ws.Cells.Item(197, 10).Value = 1400
$ws.Cells.Item(197, 11).Value = 350
$ws.Cells.Item(197, 12).Value = 400
$ws.Cells.Item(197, 13).Value = 375
$ws.Cells.Item(197, 16).Value = 188

# Row 198
$ws.Cells.Item(198, 4).Value = 44291
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 2800
$ws.Cells.Item(198, 11).Value = 450
$ws.Cells.Item(198, 12).Value = 500
$ws.Cells.Item(198, 13).Value = 475
$ws.Cells.Item(198, 16).Value = 238

# Row 199
$ws.Cells.Item(199, 4).Value = 44291
$ws.Cells.Item(199, 9).Value = "Segunda"
$ws.Cells.Item(199, 10).Value = 1460
$ws.Cells.Item(199, 11).Value = 350
$ws.Cells.Item(199, 12).Value = 400
$ws.Cells.Item(199, 13).Value = 375
$ws.Cells.Item(199, 16).Value = 188

# Row 200
$ws.Cells.Item(200, 4).Value = 44438
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 3200
$ws.Cells.Item(200, 11).Value = 450
$ws.Cells.Item(200, 12).Value = 500
$ws.Cells.Item(200, 13).Value = 475
$ws.Cells.Item(200, 16).Value = 238

# Row 201
$ws.Cells.Item(201, 4).Value = 44438
$ws.Cells.Item(201, 9).Value = "Segunda"
$ws.Cells.Item(201, 10).Value = 1540
$ws.Cells.Item(201, 11).Value = 350
$ws.Cells.Item(201, 12).Value = 400
$ws.Cells.Item(201, 13).Value = 375
$ws.Cells.Item(201, 16).Value = 188

# Row 202
$ws.Cells.Item(202, 4).Value = 44372
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 3460
$ws.Cells.Item(202, 11).Value = 500
$ws.Cells.Item(202, 12).Value = 600
$ws.Cells.Item(202, 13).Value = 550
$ws.Cells.Item(202, 16).Value = 275

# Row 203
$ws.Cells.Item(203, 4).Value = 44372
$ws.Cells.Item(203, 9).Value = "Segunda"
$ws.Cells.Item(203, 10).Value = 1600
$ws.Cells.Item(203, 11).Value = 400
$ws.Cells.Item(203, 12).Value = 450
$ws.Cells.Item(203, 13).Value = 425
$ws.Cells.Item(203, 16).Value = 212

# Row 204
$ws.Cells.Item(204, 4).Value = 44286
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 3400
$ws.Cells.Item(204, 11).Value = 450
$ws.Cells.Item(204, 12).Value = 500
$ws.Cells.Item(204, 13).Value = 475
$ws.Cells.Item(204, 16).Value = 238

# Row 205
$ws.Cells.Item(205, 4).Value = 44286
$ws.Cells.Item(205, 9).Value = "Segunda"
$ws.Cells.Item(205, 10).Value = 1600
$ws.Cells.Item(205, 11).Value = 350
$ws.Cells.Item(205, 12).Value = 400
$ws.Cells.Item(205, 13).Value = 375
$ws.Cells.Item(205, 16).Value = 188

# Row 206
$ws.Cells.Item(206, 4).Value = 44209
$ws.Cells.Item(206, 9).Value = "Primera"
$ws.Cells.Item(206, 10).Value = 3200
$ws.Cells.Item(206, 11).Value = 400
$ws.Cells.Item(206, 12).Value = 500
$ws.Cells.Item(206, 13).Value = 450
$ws.Cells.Item(206, 16).Value = 225

# Row 207
$ws.Cells.Item(207, 4).Value = 44209
$ws.Cells.Item(207, 9).Value = "Segunda"
$ws.Cells.Item(207, 10).Value = 1700
$ws.Cells.Item(207, 11).Value = 300
$ws.Cells.Item(207, 12).Value = 350
$ws.Cells.Item(207, 13).Value = 325
$ws.Cells.Item(207, 16).Value = 162

# Row 208
$ws.Cells.Item(208, 4).Value = 44356
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 3450
$ws.Cells.Item(208, 11).Value = 500
$ws.Cells.Item(208, 12).Value = 600
$ws.Cells.Item(208, 13).Value = 550
$ws.Cells.Item(208, 16).Value = 275

# Row 209
$ws.Cells.Item(209, 4).Value = 44356
$ws.Cells.Item(209, 9).Value = "Segunda"
$ws.Cells.Item(209, 10).Value = 1660
$ws.Cells.Item(209, 11).Value = 400
$ws.Cells.Item(209, 12).Value = 450
$ws.Cells.Item(209, 13).Value = 425
$ws.Cells.Item(209, 16).Value = 212

# Row 210
$ws.Cells.Item(210, 4).Value = 44160
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 2800
$ws.Cells.Item(210, 11).Value = 400
$ws.Cells.Item(210, 12).Value = 500
$ws.Cells.Item(210, 13).Value = 450
$ws.Cells.Item(210, 16).Value = 225

# Row 211
$ws.Cells.Item(211, 4).Value = 44160
$ws.Cells.Item(211, 9).Value = "Segunda"
$ws.Cells.Item(211, 10).Value = 1560
$ws.Cells.Item(211, 11).Value = 300
$ws.Cells.Item(211, 12).Value = 350
$ws.Cells.Item(211, 13).Value = 325
$ws.Cells.Item(211, 16).Value = 162

# Row 212
$ws.Cells.Item(212, 4).Value = 44351
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 3460
$ws.Cells.Item(212, 11).Value = 500
$ws.Cells.Item(212, 12).Value = 600
$ws.Cells.Item(212, 13).Value = 550
$ws.Cells.Item(212, 16).Value = 275

# Row 213
$ws.Cells.Item(213, 4).Value = 44351
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 1680
$ws.Cells.Item(213, 11).Value = 400
$ws.Cells.Item(213, 12).Value = 450
$ws.Cells.Item(213, 13).Value = 425
$ws.Cells.Item(213, 16).Value = 212

# Row 214
$ws.Cells.Item(214, 4).Value = 44365
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 3500
$ws.Cells.Item(214, 11).Value = 500
$ws.Cells.Item(214, 12).Value = 600
$ws.Cells.Item(214, 13).Value = 550
$ws.Cells.Item(214, 16).Value = 275

# Row 215
$ws.Cells.Item(215, 4).Value = 44365
$ws.Cells.Item(215, 9).Value = "Segunda"
$ws.Cells.Item(215, 10).Value = 1600
$ws.Cells.Item(215, 11).Value = 400
$ws.Cells.Item(215, 12).Value = 450
$ws.Cells.Item(215, 13).Value = 425
$ws.Cells.Item(215, 16).Value = 212

# Row 216
$ws.Cells.Item(216, 4).Value = 44162
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 2800
$ws.Cells.Item(216, 11).Value = 400
$ws.Cells.Item(216, 12).Value = 500
$ws.Cells.Item(216, 13).Value = 450
$ws.Cells.Item(216, 16).Value = 225

# Row 217
$ws.Cells.Item(217, 4).Value = 44162
$ws.Cells.Item(217, 9).Value = "Segunda"
$ws.Cells.Item(217, 10).Value = 1560
$ws.Cells.Item(217, 11).Value = 300
$ws.Cells.Item(217, 12).Value = 350
$ws.Cells.Item(217, 13).Value = 325
$ws.Cells.Item(217, 16).Value = 162

# Row 218
$ws.Cells.Item(218, 4).Value = 44410
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 3200
$ws.Cells.Item(218, 11).Value = 450
$ws.Cells.Item(218, 12).Value = 500
$ws.Cells.Item(218, 13).Value = 475
$ws.Cells.Item(218, 16).Value = 238

# Row 219
$ws.Cells.Item(219, 4).Value = 44410
$ws.Cells.Item(219, 9).Value = "Segunda"
$ws.Cells.Item(219, 10).Value = 1600
$ws.Cells.Item(219, 11).Value = 350
$ws.Cells.Item(219, 12).Value = 400
$ws.Cells.Item(219, 13).Value = 375
$ws.Cells.Item(219, 16).Value = 188

# Row 220
$ws.Cells.Item(220, 4).Value = 44244
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 3400
$ws.Cells.Item(220, 11).Value = 450
$ws.Cells.Item(220, 12).Value = 500
$ws.Cells.Item(220, 13).Value = 475
$ws.Cells.Item(220, 16).Value = 238

# Row 221
$ws.Cells.Item(221, 4).Value = 44244
$ws.Cells.Item(221, 9).Value = "Segunda"
$ws.Cells.Item(221, 10).Value = 1700
$ws.Cells.Item(221, 11).Value = 350
$ws.Cells.Item(221, 12).Value = 400
$ws.Cells.Item(221, 13).Value = 375
$ws.Cells.Item(221, 16).Value = 188

# Row 222
$ws.Cells.Item(222, 4).Value = 44176
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 2000
$ws.Cells.Item(222, 11).Value = 400
$ws.Cells.Item(222, 12).Value = 500
$ws.Cells.Item(222, 13).Value = 450
$ws.Cells.Item(222, 16).Value = 225

# Row 223
$ws.Cells.Item(223, 4).Value = 44176
$ws.Cells.Item(223, 9).Value = "Segunda"
$ws.Cells.Item(223, 10).Value = 1500
$ws.Cells.Item(223, 11).Value = 300
$ws.Cells.Item(223, 12).Value = 350
$ws.Cells.Item(223, 13).Value = 325
$ws.Cells.Item(223, 16).Value = 162

# Row 224
$ws.Cells.Item(224, 4).Value = 44239
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 3200
$ws.Cells.Item(224, 11).Value = 450
$ws.Cells.Item(224, 12).Value = 500
$ws.Cells.Item(224, 13).Value = 475
$ws.Cells.Item(224, 16).Value = 238

# Row 225
$ws.Cells.Item(225, 4).Value = 44239
$ws.Cells.Item(225, 9).Value = "Segunda"
$ws.Cells.Item(225, 10).Value = 1600
$ws.Cells.Item(225, 11).Value = 350
$ws.Cells.Item(225, 12).Value = 400
$ws.Cells.Item(225, 13).Value = 375
$ws.Cells.Item(225, 16).Value = 188

# Row 226
$ws.Cells.Item(226, 4).Value = 44358
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 3480
$ws.Cells.Item(226, 11).Value = 500
$ws.Cells.Item(226, 12).Value = 550
$ws.Cells.Item(226, 13).Value = 525
$ws.Cells.Item(226, 16).Value = 262

# Row 227
$ws.Cells.Item(227, 4).Value = 44358
$ws.Cells.Item(227, 9).Value = "Segunda"
$ws.Cells.Item(227, 10).Value = 1600
$ws.Cells.Item(227, 11).Value = 400
$ws.Cells.Item(227, 12).Value = 450
$ws.Cells.Item(227, 13).Value = 425
$ws.Cells.Item(227, 16).Value = 212

# Row 228
$ws.Cells.Item(228, 4).Value = 44211
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 3000
$ws.Cells.Item(228, 11).Value = 450
$ws.Cells.Item(228, 12).Value = 500
$ws.Cells.Item(228, 13).Value = 475
$ws.Cells.Item(228, 16).Value = 238

# Row 229
$ws.Cells.Item(229, 4).Value = 44211
$ws.Cells.Item(229, 9).Value = "Segunda"
$ws.Cells.Item(229, 10).Value = 1600
$ws.Cells.Item(229, 11).Value = 350
$ws.Cells.Item(229, 12).Value = 400
$ws.Cells.Item(229, 13).Value = 375
$ws.Cells.Item(229, 16).Value = 188

# New rows 230-231 (pushed out from the bottom of the shifted window): full row content
# Row 230
$ws.Cells.Item(230, 1).Value = 8
$ws.Cells.Item(230, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(230, 3).Value = "Coquimbo"
$ws.Cells.Item(230, 4).Value = 44323
$ws.Cells.Item(230, 5).Value = 4
$ws.Cells.Item(230, 6).Value = 100112009
$ws.Cells.Item(230, 7).Value = "Acelga"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 3440
$ws.Cells.Item(230, 11).Value = 500
$ws.Cells.Item(230, 12).Value = 600
$ws.Cells.Item(230, 13).Value = 550
$ws.Cells.Item(230, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(230, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(230, 16).Value = 275
$ws.Cells.Item(230, 17).Value = 2
$ws.Cells.Item(230, 18).Value = "Hortaliza"
$ws.Cells.Item(230, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 231
$ws.Cells.Item(231, 1).Value = 8
$ws.Cells.Item(231, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(231, 3).Value = "Coquimbo"
$ws.Cells.Item(231, 4).Value = 44323
$ws.Cells.Item(231, 5).Value = 4
$ws.Cells.Item(231, 6).Value = 100112009
$ws.Cells.Item(231, 7).Value = "Acelga"
$ws.Cells.Item(231, 8).Value = "Sin especificar"
$ws.Cells.Item(231, 9).Value = "Segunda"
$ws.Cells.Item(231, 10).Value = 1660
$ws.Cells.Item(231, 11).Value = 400
$ws.Cells.Item(231, 12).Value = 450
$ws.Cells.Item(231, 13).Value = 425
$ws.Cells.Item(231, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(231, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(231, 16).Value = 212
$ws.Cells.Item(231, 17).Value = 2
$ws.Cells.Item(231, 18).Value = "Hortaliza"
$ws.Cells.Item(231, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

